$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-17.07298043252176,1.53824973703913,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,2.522016559855211,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,1.517347111271339,-17.07298043252176,-17.07298043252176,2.785105314028501,-17.07298043252176,1.773186335543541,-17.07298043252176,2.488398886601703,-17.07298043252176),
    @(-17.07298043252176,0.8967489579773156,-17.07298043252176,-17.07298043252176,-17.07298043252176,2.060274686741856,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(2.659235031817897,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,2.032766978927989,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(3.774301001697142,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,1.552303900135581,-17.07298043252176,2.183424080751565),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,2.012266380439537,-17.07298043252176,2.447202762047524,-17.07298043252176,-17.07298043252176,-17.07298043252176,1.282805093356606),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,1.517639509202196,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,2.286968544785514,1.577766319160578),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,2.051365575942835),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,2.356412005552019,-17.07298043252176),
    @(-17.07298043252176,-1.4006379683352,4.321918154108451,-17.07298043252176,-17.07298043252176,-17.07298043252176,0.3349262706543931,0.6129523266646866,1.102359567503058,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,0.2706303794105035,0.6992632525199001,1.163524331341831,-17.07298043252176),
    @(-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,-17.07298043252176,1.780636759419155,2.045496850943455,-17.07298043252176,-17.07298043252176),
    @(-17.07298043252176,2.384979515047699,-17.07298043252176,-17.07298043252176,3.712265278382393,-17.07298043252176,2.233991290028736,2.015065008765394,-17.07298043252176,2.55934920721563),
    @(-17.07298043252176,2.758909732431904,-17.07298043252176,3.171626132641587,-17.07298043252176,3.375128813333473,2.57871072125218,-17.07298043252176,-17.07298043252176,-17.07298043252176)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
